$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.886.13"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "3.861.15"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.98%  "
$ws.Range("D7").Value = "3.857.92"
$ws.Range("E7").Value = "  +2.63%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.97"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("D15").Value = "4.501.59"
$ws.Range("D16").Value = "3.836.85"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "69.037.10"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.55"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.113"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "485.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000163"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  -3.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").Value = "4.006.55"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").Value = "3.802.81"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "437.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.97"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "2.838.85"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +12.77%  "
